$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "25.836.02"
Set-TextValue $ws.Range("E2") "  +0.22%  "
Set-TextValue $ws.Range("D3") "1.740.56"
Set-TextValue $ws.Range("E3") "  -0.41%  "
Set-TextValue $ws.Range("D4") "0.9997"
Set-TextValue $ws.Range("E4") "  -0.17%  "
Set-TextValue $ws.Range("D5") "232.17"
Set-TextValue $ws.Range("E5") "  -1.34%  "
Set-TextValue $ws.Range("D6") "1.000"
Set-TextValue $ws.Range("E6") "  -0.04%  "
Set-TextValue $ws.Range("D7") "0.5171"
Set-TextValue $ws.Range("E7") "  +1.73%  "
Set-TextValue $ws.Range("D8") "0.2817"
Set-TextValue $ws.Range("E8") "  +5.46%  "
Set-TextValue $ws.Range("D9") "39.25"
Set-TextValue $ws.Range("E9") "  -3.08%  "
Set-TextValue $ws.Range("D10") "0.06116"
Set-TextValue $ws.Range("E10") "  -1.12%  "
Set-TextValue $ws.Range("D11") "1.753.85"
Set-TextValue $ws.Range("E11") "  +0.27%  "
Set-TextValue $ws.Range("E12") "  +1.36%  "
Set-TextValue $ws.Range("D13") "15.36"
Set-TextValue $ws.Range("E13") "  -0.08%  "
Set-TextValue $ws.Range("D14") "0.6531"
Set-TextValue $ws.Range("E14") "  +5.19%  "
Set-TextValue $ws.Range("D15") "4.528"
Set-TextValue $ws.Range("E15") "  +1.27%  "
Set-TextValue $ws.Range("D16") "77.16"
Set-TextValue $ws.Range("E16") "  -0.68%  "
Set-TextValue $ws.Range("D17") "0.9992"
Set-TextValue $ws.Range("E17") "  -0.21%  "
Set-TextValue $ws.Range("D18") "0.9996"
Set-TextValue $ws.Range("E18") "  -0.06%  "
Set-TextValue $ws.Range("D19") "25.826.51"
Set-TextValue $ws.Range("E19") "  +0.08%  "
Set-TextValue $ws.Range("E20") "  -0.95%  "
Set-TextValue $ws.Range("D21") "0.000006597"
Set-TextValue $ws.Range("E21") "  -0.67%  "
Set-TextValue $ws.Range("D22") "1.977.06"
Set-TextValue $ws.Range("E22") "  +0.09%  "
Set-TextValue $ws.Range("D23") "4.129"
Set-TextValue $ws.Range("E23") "  +2.03%  "
Set-TextValue $ws.Range("D24") "8.663"
Set-TextValue $ws.Range("E24") "  +5.04%  "
Set-TextValue $ws.Range("D25") "5.151"
Set-TextValue $ws.Range("E25") "  +0.27%  "
Set-TextValue $ws.Range("D26") "139.37"
Set-TextValue $ws.Range("E26") "  +2.13%  "
Set-TextValue $ws.Range("D27") "1.512"
Set-TextValue $ws.Range("E27") "  +3.52%  "
Set-TextValue $ws.Range("D28") "15.08"
Set-TextValue $ws.Range("E28") "  +0.18%  "
Set-TextValue $ws.Range("D29") "1.812"
Set-TextValue $ws.Range("E29") "  +2.31%  "
Set-TextValue $ws.Range("D30") "102.14"
Set-TextValue $ws.Range("E30") "  -0.56%  "
Set-TextValue $ws.Range("D31") "0.08302"
Set-TextValue $ws.Range("E31") "  +1.82%  "
Set-TextValue $ws.Range("D32") "3.679"
Set-TextValue $ws.Range("E32") "  -0.41%  "
Set-TextValue $ws.Range("D33") "3.432"
Set-TextValue $ws.Range("E33") "  +1.36%  "
Set-TextValue $ws.Range("D34") "0.04500"
Set-TextValue $ws.Range("E35") "  -1.50%  "
Set-TextValue $ws.Range("D36") "0.9883"
Set-TextValue $ws.Range("E36") "  -0.58%  "
Set-TextValue $ws.Range("D37") "0.6147"
Set-TextValue $ws.Range("E37") "  +2.38%  "
Set-TextValue $ws.Range("D38") "2.646"
Set-TextValue $ws.Range("E38") "  +0.73%  "
Set-TextValue $ws.Range("D39") "0.01585"
Set-TextValue $ws.Range("E39") "  +1.98%  "
Set-TextValue $ws.Range("D40") "1.937"
Set-TextValue $ws.Range("E40") "  -0.21%  "
Set-TextValue $ws.Range("D41") "0.9993"
Set-TextValue $ws.Range("E41") "  -0.11%  "
Set-TextValue $ws.Range("D42") "100.82"
Set-TextValue $ws.Range("E42") "  -0.62%  "
Set-TextValue $ws.Range("D43") "0.3861"
Set-TextValue $ws.Range("E43") "  +1.02%  "
Set-TextValue $ws.Range("D44") "0.7257"
Set-TextValue $ws.Range("E44") "  -2.66%  "
Set-TextValue $ws.Range("D45") "4.968"
Set-TextValue $ws.Range("E45") "  +1.69%  "
Set-TextValue $ws.Range("D46") "0.05409"
Set-TextValue $ws.Range("E46") "  -1.72%  "
Set-TextValue $ws.Range("D47") "6.302"
Set-TextValue $ws.Range("E47") "  +6.48%  "
Set-TextValue $ws.Range("E48") "  +3.24%  "
Set-TextValue $ws.Range("D49") "53.26"
Set-TextValue $ws.Range("E49") "  +1.36%  "
Set-TextValue $ws.Range("D50") "7.659"
Set-TextValue $ws.Range("E50") "  +3.18%  "
Set-TextValue $ws.Range("D51") "29.89"
Set-TextValue $ws.Range("E51") "  -0.41%  "
